$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 40957.68
$ws.Range("I11").Value = 40957.68
$ws.Range("K11").Value = 40957.68
$ws.Range("M11").Value = -40817.68
$ws.Range("H41").Value = 167433.5
$ws.Range("I41").Value = 962.75
$ws.Range("J41").Value = 500375
$ws.Range("K41").Value = 962.75
$ws.Range("L41").Value = 500375
$ws.Range("M41").Value = -522.75
$ws.Range("N41").Value = -501255
$ws.Range("H42").Value = 8
$ws.Range("I42").Value = 8
$ws.Range("K42").Value = 24
$ws.Range("M42").Value = 206
$ws.Range("H43").Value = 1200
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H49").Value = 3777
$ws.Range("I49").Value = 1875
$ws.Range("J49").Value = 4468.636
$ws.Range("K49").Value = 5625
$ws.Range("L49").Value = 13405.908
$ws.Range("M49").Value = -5489
$ws.Range("N49").Value = -13677.908
$ws.Range("H63").Value = 72423.336
$ws.Range("J63").Value = 71135
$ws.Range("L63").Value = 71135
$ws.Range("N63").Value = -72383
$ws.Range("H66").Value = 72423.336
$ws.Range("J66").Value = 71135
$ws.Range("L66").Value = 213405
$ws.Range("N66").Value = -219645
$ws.Range("H69").Value = 8585.714
$ws.Range("J69").Value = 8585.714
$ws.Range("L69").Value = 25757.142
$ws.Range("N69").Value = -27505.142
$ws.Range("H72").Value = 8585.714
$ws.Range("J72").Value = 8585.714
$ws.Range("L72").Value = 77271.42600000001
$ws.Range("N72").Value = -86007.42600000001
$ws.Range("H76").Value = 5275948.5
$ws.Range("I76").Value = 4355
$ws.Range("K76").Value = 4355
$ws.Range("M76").Value = -4040
$ws.Range("H79").Value = 5275948.5
$ws.Range("I79").Value = 4355
$ws.Range("K79").Value = 4355
$ws.Range("M79").Value = -3263
$ws.Range("H86").Value = 15466858
$ws.Range("J86").Value = 18278318
$ws.Range("L86").Value = 18278318
$ws.Range("N86").Value = -18280564
$ws.Range("H89").Value = 15466858
$ws.Range("J89").Value = 18278318
$ws.Range("L89").Value = 91391590
$ws.Range("N89").Value = -91402822
$ws.Range("H137").Value = 1721.9048
$ws.Range("I137").Value = 1197.6875
$ws.Range("K137").Value = 3593.0625
$ws.Range("M137").Value = -1043.0625
$ws.Range("H138").Value = 3245.3125
$ws.Range("J138").Value = 2929.0667
$ws.Range("L138").Value = 8787.2001
$ws.Range("N138").Value = -19067.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15987.081
$ws.Range("I32").Value = 8757.467000000001
$ws.Range("J32").Value = 20916.363
$ws.Range("K32").Value = 8757.467000000001
$ws.Range("L32").Value = 20916.363
$ws.Range("M32").Value = -8470.467000000001
$ws.Range("N32").Value = -21490.363
$ws.Range("H45").Value = 3275.8333
$ws.Range("I45").Value = 3325.6667
$ws.Range("J45").Value = 3226
$ws.Range("K45").Value = 3325.6667
$ws.Range("L45").Value = 3226
$ws.Range("M45").Value = -2948.6667
$ws.Range("N45").Value = -3980
$ws.Range("H74").Value = 34489536
$ws.Range("J74").Value = 7041.1665
$ws.Range("L74").Value = 7041.1665
$ws.Range("N74").Value = -8789.166499999999
$ws.Range("H77").Value = 34489536
$ws.Range("J77").Value = 7041.1665
$ws.Range("L77").Value = 35205.8325
$ws.Range("N77").Value = -43941.8325
$ws.Range("H122").Value = 2839.7
$ws.Range("J122").Value = 4750
$ws.Range("L122").Value = 14250
$ws.Range("N122").Value = -19150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 306.66666
$ws.Range("I22").Value = 306.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 306.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -133.66666
$ws.Range("N22").ClearContents()
$ws.Range("H95").Value = 36999.6
$ws.Range("J95").Value = 36999.6
$ws.Range("L95").Value = 36999.6
$ws.Range("N95").Value = -42491.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5009.9697
$ws.Range("I31").Value = 11164.692
$ws.Range("K31").Value = 11164.692
$ws.Range("M31").Value = -10869.692
$ws.Range("H34").Value = 5009.9697
$ws.Range("I34").Value = 11164.692
$ws.Range("K34").Value = 11164.692
$ws.Range("M34").Value = -10962.692
$ws.Range("H105").Value = 896.6
$ws.Range("I105").Value = 950
$ws.Range("J105").Value = 816.5
$ws.Range("K105").Value = 950
$ws.Range("L105").Value = 816.5
$ws.Range("M105").Value = 797
$ws.Range("N105").Value = -4310.5
$ws.Range("H141").Value = 260290.58
$ws.Range("J141").Value = 260290.58
$ws.Range("L141").Value = 260290.58
$ws.Range("N141").Value = -270650.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 576.86365
$ws.Range("I2").Value = 114.5
$ws.Range("J2").Value = 1131.7
$ws.Range("K2").Value = 687
$ws.Range("L2").Value = 6790.200000000001
$ws.Range("M2").Value = -574
$ws.Range("N2").Value = -7016.200000000001
$ws.Range("H41").Value = 3555.4443
$ws.Range("I41").Value = 333.16666
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 999.4999799999999
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = -661.4999799999999
$ws.Range("N41").Value = -30676
$ws.Range("H50").Value = 799.35297
$ws.Range("I50").Value = 950.5
$ws.Range("J50").Value = 94
$ws.Range("K50").Value = 2851.5
$ws.Range("L50").Value = 282
$ws.Range("M50").Value = -2370.5
$ws.Range("N50").Value = -1244
$ws.Range("H53").Value = 799.35297
$ws.Range("I53").Value = 950.5
$ws.Range("J53").Value = 94
$ws.Range("K53").Value = 2851.5
$ws.Range("L53").Value = 282
$ws.Range("M53").Value = -2370.5
$ws.Range("N53").Value = -1244
$ws.Range("H86").Value = 166.77777
$ws.Range("I86").Value = 171.33333
$ws.Range("J86").Value = 164.5
$ws.Range("K86").Value = 513.99999
$ws.Range("L86").Value = 493.5
$ws.Range("M86").Value = 672.00001
$ws.Range("N86").Value = -2865.5
$ws.Range("H89").Value = 166.77777
$ws.Range("I89").Value = 171.33333
$ws.Range("J89").Value = 164.5
$ws.Range("K89").Value = 1541.99997
$ws.Range("L89").Value = 1480.5
$ws.Range("M89").Value = 4386.00003
$ws.Range("N89").Value = -13336.5
$ws.Range("H121").Value = 459962.7
$ws.Range("I121").Value = 842666.5
$ws.Range("J121").Value = 718.1
$ws.Range("K121").Value = 2527999.5
$ws.Range("L121").Value = 2154.3
$ws.Range("M121").Value = -2526689.5
$ws.Range("N121").Value = -4774.3
$ws.Range("H122").Value = 16667480
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 16667480
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 150007320
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -150012220
$ws.Range("H132").Value = 2505.1572
$ws.Range("I132").Value = 866.6
$ws.Range("K132").Value = 7799.400000000001
$ws.Range("M132").Value = -5269.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 333.55554
$ws.Range("I2").Value = 57.8
$ws.Range("J2").Value = 678.25
$ws.Range("K2").Value = 57.8
$ws.Range("L2").Value = 678.25
$ws.Range("M2").Value = 55.2
$ws.Range("N2").Value = -904.25
$ws.Range("H126").Value = 6782.4443
$ws.Range("I126").Value = 8585.125
$ws.Range("J126").Value = 4160.364
$ws.Range("K126").Value = 25755.375
$ws.Range("L126").Value = 12481.092
$ws.Range("M126").Value = -23285.375
$ws.Range("N126").Value = -17421.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2067.3635
$ws.Range("J46").Value = 2315.9167
$ws.Range("L46").Value = 2315.9167
$ws.Range("N46").Value = -2691.9167
$ws.Range("H136").Value = 2287.5264
$ws.Range("I136").Value = 2176.0312
$ws.Range("K136").Value = 6528.0936
$ws.Range("M136").Value = -3978.0936

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 39444.445
$ws.Range("I43").Value = 39375
$ws.Range("K43").Value = 39375
$ws.Range("M43").Value = -39226
$ws.Range("H136").Value = 2839.1794
$ws.Range("I136").Value = 2601.303
$ws.Range("J136").Value = 4147.5
$ws.Range("K136").Value = 7803.909
$ws.Range("L136").Value = 12442.5
$ws.Range("M136").Value = -5253.909
$ws.Range("N136").Value = -17542.5
